# Apply the "new version with timestamp" update to the daily-sale report.
#
# What changed (per the OOXML diff):
#   1. A new product row is inserted into the items table (sorted position
#      16, pushing the previous row 16 - "معجون سيجنال 50 مل" - down to
#      row 17 with item-number 11).
#   2. The grand-total cell (P) moves from row 17 to row 18 and its value
#      increases from 800.88 to 825.88 (+25.00 for the new item).
#   3. The footer row moves from row 18 to row 19.
#   4. The footer timestamp text changes from "11:42 AM" to "11:46 AM".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new row above the current grand-total row (row 17). ---
# This shifts: old row 17 (total) -> row 18, old row 18 (footer) -> row 19.
$ws.Rows("17:17").Insert()

# Populate the freshly-inserted row 17 by duplicating the formatting of the
# still-intact data row 16 (same per-column styles as every other item row).
$ws.Range("A16:Q16").Copy($ws.Range("A17:Q17"))
$ws.Rows("17:17").RowHeight = 25.5

# Re-merge the cells of the new row exactly like every other item row.
$ws.Range("A17:B17").Merge()
$ws.Range("C17:G17").Merge()
$ws.Range("H17:K17").Merge()
$ws.Range("L17:M17").Merge()
$ws.Range("N17:O17").Merge()

# Row 17 now holds a copy of what used to be row 16 ("معجون سيجنال 50 مل").
# Renumber it as item 11 (row 16 keeps item number 10).
$ws.Range("A17").Value = 11
$ws.Range("C17").Value = "معجون سيجنال 50 مل"
$ws.Range("H17").Value = "8:0"
$ws.Range("L17").Value = "0"
$ws.Range("N17").Value = "35.00"
$ws.Range("P17").Value = "35.0000"
$ws.Range("Q17").Value = "1:0"

# --- 2. Row 16 becomes the newly-added product. ---
$ws.Range("C16").Value = "معجون حلاقه 55555"
$ws.Range("H16").Value = "12:0"
$ws.Range("N16").Value = "25.00"
$ws.Range("P16").Value = "25.0000"

# --- 3. Update the grand total (now on row 18) for the new item's value. ---
$ws.Range("P18").Value = 825.88
$ws.Rows("18:18").RowHeight = 24.75

# --- 4. Update the footer timestamp (now on row 19). ---
$ws.Range("A19").Value = "Friday, 5 September, 2025 11:46 AM"
